# Generate Report for Handoff
# Applies the localization-status.xlsx refresh: new source file uuid
# (0a24f56a... -> 1245ae1e...), new handoff package hash
# (d72723ae... -> 8811a66b5...), refreshed timestamps, and clears the
# already-consumed "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the zh-cn and de-de sheets
# because the new handback has not happened yet.

$wb = $excel.ActiveWorkbook

$oldGuid = "0a24f56a-195f-4297-a987-da646670bcb0"
$newGuid = "1245ae1e-ea73-4968-94a4-9be32f154c0f"
$newHash = "8811a66b534f636e99153072e93851906514d686"

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws1.Range("A2").Value = "$newGuid.md"
$ws1.Range("B2").Value = "e2e\$newGuid.md"
$ws1.Range("G2").Value = "2016-09-01 19:10:19"

$ov_addr = $null
foreach ($hl in $ws1.Hyperlinks) {
    $ov_addr = $hl.Address
}
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $ov_addr, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md") | Out-Null

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$a2addr = $null
foreach ($hl in $ws2.Hyperlinks) {
    $r = $hl.Range.Address()
    if ($r -eq '$A$2') {
        $a2addr = $hl.Address
    }
}
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $a2addr, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null

$ws2.Range("A2").Value = "$newGuid.md"
$ws2.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-01 19:10:01"
$ws2.Range("I2").ClearContents()
$ws2.Range("I2").Style = "Normal"
$ws2.Range("J2").ClearContents()
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

$ws2.Columns.Item(9).ColumnWidth = 17.83
$ws2.Columns.Item(10).ColumnWidth = 20.83

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$a2addr3 = $null
foreach ($hl in $ws3.Hyperlinks) {
    $r = $hl.Range.Address()
    if ($r -eq '$A$2') {
        $a2addr3 = $hl.Address
    }
}
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $a2addr3, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null

$ws3.Range("A2").Value = "$newGuid.md"
$ws3.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-01 19:10:19"
$ws3.Range("I2").ClearContents()
$ws3.Range("I2").Style = "Normal"
$ws3.Range("J2").ClearContents()
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$ws3.Columns.Item(9).ColumnWidth = 17.83
$ws3.Columns.Item(10).ColumnWidth = 20.83
